$p = $ppt.ActivePresentation

# --- Edit 1: Slide 21 ("Declarative Programming" title) -----------------
# Originally two runs: "Declarative " + "Programming".
# Target: a single run "Declarative Programming" that keeps the SECOND
# run's rPr (lang="en-US" dirty="0"). We delete the text of the first
# run (collapsing it) and then prepend the removed text onto the
# remaining run via InsertBefore, which makes the combined text inherit
# the still-alive (second) run's formatting.
$slide21 = $p.Slides.Item(21)
$titleShape = $slide21.Shapes.Item(1)
$titleTextRange = $titleShape.TextFrame.TextRange

$firstWordLen = "Declarative ".Length
$firstRunRange = $titleTextRange.Characters(1, $firstWordLen)
$firstRunRange.Text = ""

$remainingRange = $titleTextRange.Characters(1, "Programming".Length)
$remainingRange.InsertBefore("Declarative ") | Out-Null

# --- Edit 2: Slide 8 (code sample "lst += [ 5 ]" -> "lst += [ i*i ]") ---
$slide8 = $p.Slides.Item(8)
$codeShape = $slide8.Shapes.Item(2)
$codeTextRange = $codeShape.TextFrame.TextRange

$fullText = $codeTextRange.Text
$anchor = $fullText.IndexOf("lst += [ 5 ]")
$fiveStart = $anchor + 1 + "lst += [ ".Length

# Turn the literal "5" into the first "i".
$codeTextRange.Characters($fiveStart, 1).Text = "i"

# Insert "*i" right after the first "i" (covers the "*" and the second "i").
$codeTextRange.Characters($fiveStart, 1).InsertAfter("*i") | Out-Null

# The newly inserted "*i" currently shares one run because it carries
# identical formatting to its neighbours. Nudge the "*" character's
# formatting (off and back on) so it keeps its own run, which in turn
# forces the trailing "i" to split off into its own run as well - giving
# the five distinct runs ("lst", " += [ ", "i", "*", "i", " ]") that the
# target markup uses.
$starRange = $codeTextRange.Characters($fiveStart + 1, 1)
$starRange.Font.Bold = $true
$starRange.Font.Bold = $false
